$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 120

$ws.Range("H3").Value = 200
$ws.Range("I3").Value = 240

$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 360

$ws.Range("H5").Value = 400
$ws.Range("I5").Value = 480

$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 600

$null = $ws.Range("I7").Select()
